# Append the Aug. 8, 2022 data row (data updated on Aug.09) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 26

# Date column keeps the same "M/D/YYYY" number format as the rows above it.
$ws.Cells.Item($newRow, 1).Value = 44781
$ws.Cells.Item($newRow, 1).NumberFormat = "M/D/YYYY"

# New day's counts.
$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 0
$ws.Cells.Item($newRow, 4).Value = 0

# Match the author's final selection on the newly entered row.
[void]$ws.Range("B26:D26").Select()
